$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 783060.1
$ws.Range("I15").Value = 783060.1
$ws.Range("K15").Value = 2349180.3
$ws.Range("M15").Value = -2349011.3
$ws.Range("H100").Value = 1405.2
$ws.Range("I100").Value = 1506.25
$ws.Range("K100").Value = 1506.25
$ws.Range("M100").Value = -965.25
$ws.Range("H101").Value = 399.75
$ws.Range("J101").Value = 259
$ws.Range("L101").Value = 777
$ws.Range("N101").Value = -4021
$ws.Range("H115").Value = 519.6667
$ws.Range("I115").Value = 519.6667
$ws.Range("K115").Value = 1559.0001
$ws.Range("M115").Value = 7.999900000000025
$ws.Range("H118").Value = 1279.8889
$ws.Range("I118").Value = 1279.8889
$ws.Range("K118").Value = 3839.6667
$ws.Range("M118").Value = -2182.6667
$ws.Range("H127").Value = 2227.1667
$ws.Range("I127").Value = 2032.6
$ws.Range("K127").Value = 6097.799999999999
$ws.Range("M127").Value = -1137.799999999999
$ws.Range("H129").Value = 4680.2
$ws.Range("J129").Value = 2076.25
$ws.Range("L129").Value = 6228.75
$ws.Range("N129").Value = -16228.75
$ws.Range("H132").Value = 2374.6553
$ws.Range("I132").Value = 2295.2222
$ws.Range("K132").Value = 6885.6666
$ws.Range("M132").Value = -4355.6666
$ws.Range("H137").Value = 8196.764999999999
$ws.Range("I137").Value = 6612.1665
$ws.Range("J137").Value = 11999.8
$ws.Range("K137").Value = 19836.4995
$ws.Range("L137").Value = 35999.39999999999
$ws.Range("M137").Value = -17286.4995
$ws.Range("N137").Value = -41099.39999999999
$ws.Range("H138").Value = 3142.279
$ws.Range("J138").Value = 4496
$ws.Range("L138").Value = 13488
$ws.Range("N138").Value = -23768

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 459.1111
$ws.Range("I5").Value = 300.3913
$ws.Range("K5").Value = 300.3913
$ws.Range("M5").Value = -188.3913
$ws.Range("H32").Value = 7527.759
$ws.Range("I32").Value = 7060.826
$ws.Range("J32").Value = 9317.666999999999
$ws.Range("K32").Value = 7060.826
$ws.Range("L32").Value = 9317.666999999999
$ws.Range("M32").Value = -6773.826
$ws.Range("N32").Value = -9891.666999999999
$ws.Range("H63").Value = 2008.6
$ws.Range("I63").Value = 1975
$ws.Range("J63").Value = 2031
$ws.Range("K63").Value = 1975
$ws.Range("L63").Value = 2031
$ws.Range("M63").Value = -1289
$ws.Range("N63").Value = -3403
$ws.Range("H66").Value = 2008.6
$ws.Range("I66").Value = 1975
$ws.Range("J66").Value = 2031
$ws.Range("K66").Value = 9875
$ws.Range("L66").Value = 10155
$ws.Range("M66").Value = -6443
$ws.Range("N66").Value = -17019
$ws.Range("H102").Value = 3334.353
$ws.Range("I102").Value = 2788.1428
$ws.Range("J102").Value = 5883.3335
$ws.Range("K102").Value = 2788.1428
$ws.Range("L102").Value = 5883.3335
$ws.Range("M102").Value = -1166.1428
$ws.Range("N102").Value = -9127.333500000001
$ws.Range("H132").Value = 3315.0952
$ws.Range("J132").Value = 9299.799999999999
$ws.Range("L132").Value = 27899.4
$ws.Range("N132").Value = -32959.39999999999
$ws.Range("H137").Value = 85000
$ws.Range("J137").Value = 85000
$ws.Range("L137").Value = 85000
$ws.Range("N137").Value = -95200

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 459.1111
$ws.Range("I4").Value = 300.3913
$ws.Range("K4").Value = 300.3913
$ws.Range("M4").Value = -185.3913
$ws.Range("H22").Value = 689.25
$ws.Range("I22").Value = 699
$ws.Range("J22").Value = 582
$ws.Range("K22").Value = 699
$ws.Range("L22").Value = 582
$ws.Range("M22").Value = -526
$ws.Range("N22").Value = -928
$ws.Range("H134").Value = 3133.6875
$ws.Range("I134").Value = 1549.6154
$ws.Range("J134").Value = 9998
$ws.Range("K134").Value = 4648.8462
$ws.Range("L134").Value = 29994
$ws.Range("M134").Value = -2113.8462
$ws.Range("N134").Value = -35064

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 797
$ws.Range("J11").Value = 692.5
$ws.Range("L11").Value = 692.5
$ws.Range("N11").Value = -972.5
$ws.Range("H31").Value = 6542.4062
$ws.Range("J31").Value = 21373.285
$ws.Range("L31").Value = 21373.285
$ws.Range("N31").Value = -21963.285
$ws.Range("H34").Value = 6542.4062
$ws.Range("J34").Value = 21373.285
$ws.Range("L34").Value = 21373.285
$ws.Range("N34").Value = -21777.285
$ws.Range("H133").Value = 55000
$ws.Range("I133").Value = 50000
$ws.Range("J133").Value = 60000
$ws.Range("K133").Value = 50000
$ws.Range("L133").Value = 60000
$ws.Range("M133").Value = -47470
$ws.Range("N133").Value = -65060
$ws.Range("H135").Value = 103423
$ws.Range("I135").Value = 60709
$ws.Range("J135").Value = 124780
$ws.Range("K135").Value = 60709
$ws.Range("L135").Value = 124780
$ws.Range("M135").Value = -55639
$ws.Range("N135").Value = -134920
$ws.Range("H140").Value = 80750
$ws.Range("J140").Value = 80750
$ws.Range("L140").Value = 80750
$ws.Range("N140").Value = -91110

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 81.65000000000001
$ws.Range("I2").Value = 113.416664
$ws.Range("J2").Value = 34
$ws.Range("K2").Value = 680.499984
$ws.Range("L2").Value = 204
$ws.Range("M2").Value = -567.499984
$ws.Range("N2").Value = -430
$ws.Range("H38").Value = 43.857143
$ws.Range("I38").Value = 40
$ws.Range("J38").Value = 67
$ws.Range("K38").Value = 120
$ws.Range("L38").Value = 201
$ws.Range("M38").Value = 227
$ws.Range("N38").Value = -895
$ws.Range("H44").Value = 786.125
$ws.Range("I44").Value = 750.5
$ws.Range("J44").Value = 893
$ws.Range("K44").Value = 2251.5
$ws.Range("L44").Value = 2679
$ws.Range("M44").Value = -1853.5
$ws.Range("N44").Value = -3475
$ws.Range("H109").Value = 434
$ws.Range("I109").Value = 237.27272
$ws.Range("J109").Value = 975
$ws.Range("K109").Value = 711.81816
$ws.Range("L109").Value = 2925
$ws.Range("M109").Value = 328.18184
$ws.Range("N109").Value = -5005
$ws.Range("H122").Value = 810.94116
$ws.Range("I122").Value = 769.8889
$ws.Range("J122").Value = 857.125
$ws.Range("K122").Value = 6929.0001
$ws.Range("L122").Value = 7714.125
$ws.Range("M122").Value = -4479.0001
$ws.Range("N122").Value = -12614.125
$ws.Range("H128").Value = 524998.75
$ws.Range("I128").Value = 524998.75
$ws.Range("K128").Value = 1574996.25
$ws.Range("M128").Value = -1570016.25
$ws.Range("H131").Value = 1138171.2
$ws.Range("I131").Value = 664.6667
$ws.Range("K131").Value = 1994.0001
$ws.Range("M131").Value = 3045.9999
$ws.Range("H141").Value = 4600.8335
$ws.Range("I141").Value = 4600.8335
$ws.Range("K141").Value = 13802.5005
$ws.Range("M141").Value = -8622.500499999998

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 6447.3335
$ws.Range("I9").Value = 7450
$ws.Range("J9").Value = 4442
$ws.Range("K9").Value = 7450
$ws.Range("L9").Value = 4442
$ws.Range("M9").Value = -7280
$ws.Range("N9").Value = -4782
$ws.Range("H132").Value = 3316.625
$ws.Range("I132").Value = 2390.4285
$ws.Range("J132").Value = 9800
$ws.Range("K132").Value = 7171.2855
$ws.Range("L132").Value = 29400
$ws.Range("M132").Value = -4641.2855
$ws.Range("N132").Value = -34460

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7939.769
$ws.Range("I40").Value = 6110.727
$ws.Range("K40").Value = 6110.727
$ws.Range("M40").Value = -5974.727
$ws.Range("H122").Value = 5504.067
$ws.Range("I122").Value = 5361.5
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 16084.5
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -13634.5
$ws.Range("N122").Value = -27400
$ws.Range("H136").Value = 6467
$ws.Range("I136").Value = 2131.8125
$ws.Range("J136").Value = 7974.891
$ws.Range("K136").Value = 6395.4375
$ws.Range("L136").Value = 23924.673
$ws.Range("M136").Value = -3845.4375
$ws.Range("N136").Value = -29024.673

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6859.6
$ws.Range("I132").Value = 4074.5
$ws.Range("J132").Value = 18000
$ws.Range("K132").Value = 12223.5
$ws.Range("L132").Value = 54000
$ws.Range("M132").Value = -9693.5
$ws.Range("N132").Value = -59060
$ws.Range("H136").Value = 3105.537
$ws.Range("I136").Value = 2647.7273
$ws.Range("J136").Value = 5119.9
$ws.Range("K136").Value = 7943.1819
$ws.Range("L136").Value = 15359.7
$ws.Range("M136").Value = -5393.1819
$ws.Range("N136").Value = -20459.7
